# Overview_extinction sheet: add "Gum"-rounded (significant-digit rounded)
# columns AN:AQ mirroring AI:AL ("C P Tot", "C P > 2", "C P 2 - 0.5",
# "C P 0.5 - 0.08"), plus a small "Sign Stellen" (significant digits) helper
# cell, then clear the stray old "Conf" header that is no longer used.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview_extinction")

# --- Header row (row 1): mirror AI1:AL1 into AN1:AQ1, keep AR1 text ---
$ws.Range("AN1:AQ1").Value = $ws.Range("AI1:AL1").Value

# --- Row 54: label + significant-digit count used by the rounding formulas ---
$ws.Range("AN54").Value = "Sign Stellen"
$ws.Range("AO54").Value = 2

# --- Rows 2-26: rounded "mapped" values, referencing the digit-count rows 28-52 ---
for ($r = 2; $r -le 26; $r++) {
    $src = $r + 26
    $ws.Range("AN$r").Formula = "=ROUNDUP(AI$r,-INT(LOG(AI$r))+(LEN(AN$src)-LEN(SUBSTITUTE(AN$src,`"0`",))+1)-1)"
    $ws.Range("AO$r").Formula = "=ROUNDUP(AJ$r,-INT(LOG(AJ$r))+(LEN(AO$src)-LEN(SUBSTITUTE(AO$src,`"0`",))+1)-1)"
    $ws.Range("AP$r").Formula = "=ROUNDUP(AK$r,-INT(LOG(AK$r))+(LEN(AP$src)-LEN(SUBSTITUTE(AP$src,`"0`",))+1)-1)"
    $ws.Range("AQ$r").Formula = "=ROUNDUP(AL$r,-INT(LOG(AL$r))+(LEN(AQ$src)-LEN(SUBSTITUTE(AQ$src,`"0`",))+1)-1)"
}

# --- Rows 28-52: significant-digit counts derived from the stdev rows ---
for ($r = 28; $r -le 52; $r++) {
    $ws.Range("AN$r").Formula = "=IFERROR(ROUNDUP(AI$r,-INT(LOG(AI$r))+(`$AO`$54-1)),0)"
    $ws.Range("AO$r").Formula = "=IFERROR(ROUNDUP(AJ$r,-INT(LOG(AJ$r))+(`$AO`$54-1)),0)"
    $ws.Range("AP$r").Formula = "=IFERROR(ROUNDUP(AK$r,-INT(LOG(AK$r))+(`$AO`$54-1)),0)"
    $ws.Range("AQ$r").Formula = "=IFERROR(ROUNDUP(AL$r,-INT(LOG(AL$r))+(`$AO`$54-1)),0)"
}

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("AQ2:AQ26").Select
